$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 88.8
$ws.Range("I4").Value = 88.8
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 88.8
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 25.2
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9034.270500000001
$ws.Range("I43").Value = 8157.8423
$ws.Range("J43").Value = 9959.388999999999
$ws.Range("K43").Value = 8157.8423
$ws.Range("L43").Value = 9959.388999999999
$ws.Range("M43").Value = -8088.8423
$ws.Range("N43").Value = -10097.389

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2168.524
$ws.Range("I80").Value = 255.77777
$ws.Range("J80").Value = 3603.0833
$ws.Range("K80").Value = 767.33331
$ws.Range("L80").Value = 10809.2499
$ws.Range("M80").Value = 230.66669
$ws.Range("N80").Value = -12805.2499

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2168.524
$ws.Range("I83").Value = 255.77777
$ws.Range("J83").Value = 3603.0833
$ws.Range("K83").Value = 2301.99993
$ws.Range("L83").Value = 32427.7497
$ws.Range("M83").Value = 2690.00007
$ws.Range("N83").Value = -42411.7497

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 11755.421
$ws.Range("I106").Value = 3928.889
$ws.Range("K106").Value = 3928.889
$ws.Range("M106").Value = -3297.889

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1400.2142
$ws.Range("J125").Value = 1189.25
$ws.Range("L125").Value = 10703.25
$ws.Range("N125").Value = -15623.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1106.6567
$ws.Range("I132").Value = 1084.5
$ws.Range("K132").Value = 3253.5
$ws.Range("M132").Value = -723.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3506.6667
$ws.Range("I88").Value = 6333.3335
$ws.Range("J88").Value = 2800
$ws.Range("K88").Value = 6333.3335
$ws.Range("L88").Value = 2800
$ws.Range("M88").Value = -5927.3335
$ws.Range("N88").Value = -3612

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3506.6667
$ws.Range("I91").Value = 6333.3335
$ws.Range("J91").Value = 2800
$ws.Range("K91").Value = 6333.3335
$ws.Range("L91").Value = 2800
$ws.Range("M91").Value = -4929.3335
$ws.Range("N91").Value = -5608

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2208.0908
$ws.Range("I102").Value = 2208.0908
$ws.Range("K102").Value = 2208.0908
$ws.Range("M102").Value = -586.0907999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 111112490
$ws.Range("I122").Value = 1550.875
$ws.Range("K122").Value = 4652.625
$ws.Range("M122").Value = -2202.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2146.4426
$ws.Range("I132").Value = 1380.3214
$ws.Range("K132").Value = 4140.9642
$ws.Range("M132").Value = -1610.9642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33643.383
$ws.Range("I31").Value = 1993.0834
$ws.Range("K31").Value = 1993.0834
$ws.Range("M31").Value = -1698.0834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 33643.383
$ws.Range("I34").Value = 1993.0834
$ws.Range("K34").Value = 1993.0834
$ws.Range("M34").Value = -1791.0834

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 5121.2
$ws.Range("I105").Value = 4891.909
$ws.Range("K105").Value = 4891.909
$ws.Range("M105").Value = -3144.909

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1028
$ws.Range("I107").Value = 924.0909
$ws.Range("K107").Value = 924.0909
$ws.Range("M107").Value = 995.9091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9957
$ws.Range("I132").Value = 900
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 2700
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -170
$ws.Range("N132").Value = -62102

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3297.3333
$ws.Range("I132").Value = 1422.3636
$ws.Range("K132").Value = 12801.2724
$ws.Range("M132").Value = -10271.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 8342499.5
$ws.Range("I11").Value = 8755000
$ws.Range("K11").Value = 8755000
$ws.Range("M11").Value = -8754861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11116779
$ws.Range("I102").Value = 16670169
$ws.Range("K102").Value = 16670169
$ws.Range("M102").Value = -16668547

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 67842.664
$ws.Range("J106").Value = 67842.664
$ws.Range("L106").Value = 67842.664
$ws.Range("N106").Value = -70366.664

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3373
$ws.Range("I122").Value = 1498.75
$ws.Range("K122").Value = 4496.25
$ws.Range("M122").Value = -2046.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6141
$ws.Range("I126").Value = 4194.6
$ws.Range("J126").Value = 11007
$ws.Range("K126").Value = 12583.8
$ws.Range("L126").Value = 33021
$ws.Range("M126").Value = -10113.8
$ws.Range("N126").Value = -37961

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 939644
$ws.Range("I132").Value = 1458617.5
$ws.Range("J132").Value = 5491.8
$ws.Range("K132").Value = 4375852.5
$ws.Range("L132").Value = 16475.4
$ws.Range("M132").Value = -4373322.5
$ws.Range("N132").Value = -21535.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2875.1
$ws.Range("I16").Value = 2300
$ws.Range("J16").Value = 5175.5
$ws.Range("K16").Value = 2300
$ws.Range("L16").Value = 5175.5
$ws.Range("M16").Value = -2130
$ws.Range("N16").Value = -5515.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7205.4
$ws.Range("I40").Value = 6362.5713
$ws.Range("K40").Value = 6362.5713
$ws.Range("M40").Value = -6226.5713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3664.9285
$ws.Range("I82").Value = 5692.375
$ws.Range("J82").Value = 2853.95
$ws.Range("K82").Value = 5692.375
$ws.Range("L82").Value = 2853.95
$ws.Range("M82").Value = -5331.375
$ws.Range("N82").Value = -3575.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3664.9285
$ws.Range("I85").Value = 5692.375
$ws.Range("J85").Value = 2853.95
$ws.Range("K85").Value = 5692.375
$ws.Range("L85").Value = 2853.95
$ws.Range("M85").Value = -4444.375
$ws.Range("N85").Value = -5349.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2000
$ws.Range("J93").Value = 2000
$ws.Range("L93").Value = 2000
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 9494.833000000001
$ws.Range("I100").Value = 6730.2104
$ws.Range("J100").Value = 12584.706
$ws.Range("K100").Value = 6730.2104
$ws.Range("L100").Value = 12584.706
$ws.Range("M100").Value = -6189.2104
$ws.Range("N100").Value = -13666.706

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11952
$ws.Range("I122").Value = 4899
$ws.Range("J122").Value = 19005
$ws.Range("K122").Value = 14697
$ws.Range("L122").Value = 57015
$ws.Range("M122").Value = -12247
$ws.Range("N122").Value = -61915

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7274.1113
$ws.Range("I132").Value = 3899.4375
$ws.Range("J132").Value = 12182.728
$ws.Range("K132").Value = 11698.3125
$ws.Range("L132").Value = 36548.18399999999
$ws.Range("M132").Value = -9168.3125
$ws.Range("N132").Value = -41608.18399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8637.177
$ws.Range("I136").Value = 3715.5557
$ws.Range("J136").Value = 14174
$ws.Range("K136").Value = 11146.6671
$ws.Range("L136").Value = 42522
$ws.Range("M136").Value = -8596.667099999999
$ws.Range("N136").Value = -47622

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 76248
$ws.Range("J139").Value = 76248
$ws.Range("L139").Value = 76248
$ws.Range("N139").Value = -86528

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 9333.333000000001
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 9333.333000000001
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 9333.333000000001
$ws.Range("N4").Value = -9559.333000000001
$ws.Range("M4").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3109.4443
$ws.Range("I81").Value = 2174.5881
$ws.Range("J81").Value = 19002
$ws.Range("K81").Value = 4349.1762
$ws.Range("L81").Value = 38004
$ws.Range("M81").Value = -3288.1762
$ws.Range("N81").Value = -40126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3109.4443
$ws.Range("I84").Value = 2174.5881
$ws.Range("J84").Value = 19002
$ws.Range("K84").Value = 21745.881
$ws.Range("L84").Value = 190020
$ws.Range("M84").Value = -16441.881
$ws.Range("N84").Value = -200628

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 48000
$ws.Range("J95").Value = 48000
$ws.Range("L95").Value = 48000
$ws.Range("N95").Value = -53492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1874.5
$ws.Range("J96").Value = 2089.4
$ws.Range("L96").Value = 2089.4
$ws.Range("N96").Value = -4835.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 50370
$ws.Range("J104").Value = 50370
$ws.Range("L104").Value = 50370
$ws.Range("N104").Value = -57358

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4827
$ws.Range("I122").Value = 1991.4
$ws.Range("J122").Value = 19005
$ws.Range("K122").Value = 5974.200000000001
$ws.Range("L122").Value = 57015
$ws.Range("M122").Value = -3524.200000000001
$ws.Range("N122").Value = -61915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1743.1111
$ws.Range("I136").Value = 1249.9143
$ws.Range("J136").Value = 19005
$ws.Range("K136").Value = 3749.7429
$ws.Range("L136").Value = 57015
$ws.Range("M136").Value = -1199.7429
$ws.Range("N136").Value = -62115
